$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 8 with a textual date value "23-02-2019" and trailing
# numeric stats matching the existing table layout.
$ws.Range("A8").Value = "23-02-2019"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 346
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
